# Update the Training Dashboard sheet with new progress as of 04-Nov-2025:
#  - Column H (PERIOD TO EXPIRE) decreases by 1 for each data row (one day closer to expiry)
#  - Column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 19; $row++) {
    $cellH = $ws.Cells.Item($row, 8)   # column H
    $cellI = $ws.Cells.Item($row, 9)   # column I

    if ($cellH.Value2 -ne $null) {
        $cellH.Value2 = $cellH.Value2 - 1
    }

    if ($cellI.Value2 -eq "03-Nov-2025") {
        # Assign via Formula with a leading apostrophe so the date-like
        # string "04-Nov-2025" is stored as literal text (matching the
        # original inline-string cell) rather than being auto-converted
        # into a date serial number by Excel's input parser.
        $cellI.Formula = "'04-Nov-2025"
    }
}

$wb.Save()
